$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Job# values for B2:B11 (each original value shifted by +30,
# e.g. 32255216 -> 32255246 ... 32255225 -> 32255255)
$newValues = @(
    "32255246",
    "32255247",
    "32255248",
    "32255249",
    "32255250",
    "32255251",
    "32255252",
    "32255253",
    "32255254",
    "32255255"
)

$tempCell = $ws.Range("Z1")

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $target = $ws.Range("B$row")

    # Build the value as a formula-derived text result first (via a scratch
    # cell) and paste only the value into the target cell. This keeps the
    # cell stored as a plain shared-string ("t=s") without Excel adding a
    # "quote prefix" style the way a direct apostrophe-prefixed literal
    # would, matching how the original numeric-looking Job# strings are
    # stored in the workbook.
    $tempCell.Formula = '="' + $newValues[$i] + '"'
    $tempCell.Copy()
    $target.PasteSpecial(-4163)
}

$tempCell.Clear()
$excel.CutCopyMode = $false
